$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rewrite existing headers and add new ones ---
$ws.Range("A1").Value2 = "Graph"
$ws.Range("B1").Value2 = "('Betweenness', 'node1')"
$ws.Range("C1").Value2 = "('Betweenness', 'value1')"
$ws.Range("D1").Value2 = "('Betweenness', 'node2')"
$ws.Range("E1").Value2 = "('Betweenness', 'value2')"
$ws.Range("F1").Value2 = "('Betweenness', 'node3')"
$ws.Range("G1").Value2 = "('Betweenness', 'value3')"
$ws.Range("H1").Value2 = "('Closeness', 'node1')"
$ws.Range("I1").Value2 = "('Closeness', 'value1')"
$ws.Range("J1").Value2 = "('Closeness', 'node2')"
$ws.Range("K1").Value2 = "('Closeness', 'value2')"
$ws.Range("L1").Value2 = "('Closeness', 'node3')"
$ws.Range("M1").Value2 = "('Closeness', 'value3')"
$ws.Range("N1").Value2 = "('Degree', 'node1')"
$ws.Range("O1").Value2 = "('Degree', 'value1')"
$ws.Range("P1").Value2 = "('Degree', 'node2')"
$ws.Range("Q1").Value2 = "('Degree', 'value2')"
$ws.Range("R1").Value2 = "('Degree', 'node3')"
$ws.Range("S1").Value2 = "('Degree', 'value3')"

# Apply the same header format (bold, centered, thin border) used by A1:G1
# to the newly added header cells H1:S1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("H1:S1").PasteSpecial(-4122) | Out-Null

# --- Row 2 (Astronomy) ---
$ws.Range("A2").Value2 = "Astronomy"
$ws.Range("B2").Value2 = "태양계"
$ws.Range("C2").Value2 = 0.1983368870485608
$ws.Range("D2").Value2 = "케플러법칙"
$ws.Range("E2").Value2 = 0.1478279317385609
$ws.Range("F2").Value2 = "행성"
$ws.Range("G2").Value2 = 0.1410976472855235
$ws.Range("H2").Value2 = "행성"
$ws.Range("I2").Value2 = 0.3500268749776042
$ws.Range("J2").Value2 = "태양계"
$ws.Range("K2").Value2 = 0.3492884216548666
$ws.Range("L2").Value2 = "케플러법칙"
$ws.Range("M2").Value2 = 0.3083104503992677
$ws.Range("N2").Value2 = "태양계"
$ws.Range("O2").Value2 = 82
$ws.Range("P2").Value2 = "행성"
$ws.Range("Q2").Value2 = 80
$ws.Range("R2").Value2 = "케플러법칙"
$ws.Range("S2").Value2 = 70

# --- Row 3 (Sampling) ---
$ws.Range("A3").Value2 = "Sampling"
$ws.Range("B3").Value2 = "확률적표본추출"
$ws.Range("C3").Value2 = 0.2296290410994625
$ws.Range("D3").Value2 = "비확률적표본추출"
$ws.Range("E3").Value2 = 0.202396998207084
$ws.Range("F3").Value2 = "표본추출"
$ws.Range("G3").Value2 = 0.1656861221582649
$ws.Range("H3").Value2 = "비확률적표본추출"
$ws.Range("I3").Value2 = 0.3466749191940032
$ws.Range("J3").Value2 = "확률적표본추출"
$ws.Range("K3").Value2 = 0.3453567636077142
$ws.Range("L3").Value2 = "표본추출"
$ws.Range("M3").Value2 = 0.3186976450134345
$ws.Range("N3").Value2 = "확률적표본추출"
$ws.Range("O3").Value2 = 98
$ws.Range("P3").Value2 = "비확률적표본추출"
$ws.Range("Q3").Value2 = 88
$ws.Range("R3").Value2 = "표본추출"
$ws.Range("S3").Value2 = 82

# --- Row 4 (Database) ---
$ws.Range("A4").Value2 = "Database"
$ws.Range("B4").Value2 = "RDBMS"
$ws.Range("C4").Value2 = 0.2291077994660412
$ws.Range("D4").Value2 = "데이터베이스"
$ws.Range("E4").Value2 = 0.2008070968266106
$ws.Range("F4").Value2 = "non-RDBMS"
$ws.Range("G4").Value2 = 0.1603612529923504
$ws.Range("H4").Value2 = "RDBMS"
$ws.Range("I4").Value2 = 0.3687791566020688
$ws.Range("J4").Value2 = "데이터베이스"
$ws.Range("K4").Value2 = 0.3585610943284599
$ws.Range("L4").Value2 = "non-RDBMS"
$ws.Range("M4").Value2 = 0.3480409644215124
$ws.Range("N4").Value2 = "RDBMS"
$ws.Range("O4").Value2 = 141
$ws.Range("P4").Value2 = "데이터베이스"
$ws.Range("Q4").Value2 = 133
$ws.Range("R4").Value2 = "non-RDBMS"
$ws.Range("S4").Value2 = 120
